$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCaseMatrix")

# Update "Actual result if unexpected" (E) and "Result" (F) for the
# "maximal edge case" row (row 5)
$ws.Range("E5").Value = "*** Start of Concatenating Strings Demo ***`nType the 1st string (q - to quit):`nQbaKozj6CjkH1ATLPQ6rglly9riWzt3nRbrCrrmN5DfMxC8PhwgZSn5vDKdkMPDvBzM2fOJ2lBZum8YbINYUGGlALWwJGGlhK5VQoJ2Xadht4K7hxwb8ChkIbSHdITCdsdHLXm3LL9lwCsWICIb19TNwjrCPwervA43DOAD7KMlnu5lPnWJ4ca9ua4vTUgQ8EwQAKHcg`nType the 2nd string:`nConcatenated string is 'QbaKozj6CjkH1ATLPQ6rglly9riWzt3nRbrCrrmN5DfMxC8PhwgZSn5vDKdkMPDvBzM2fOJ2lBZum8YbINYUGGlALWwJGGlhK5QoJ2Xadht4K7hxwb8ChkIbSHdITCdsdHLXm3LL9lwCsWICIb19TNwjrCPwervA43DOAD7KMlnu5lPnWJ4ca9ua4vTUgQ8EwQAK'`nType the 1st string (q - to quit):`nType the 2nd string:"
$ws.Range("F5").Value = "FAIL"

# Update "Actual result if unexpected" (E), "Result" (F), and
# "Comments" (G) for the "minimal edge case" row (row 15)
$ws.Range("E15").Value = "*** Start of Concatenating Strings Demo ***`nType the 1st string (q - to quit):`nZyBk27e3yEyMQ0WvL8t6KIJauYhrPmWuyRC2bnNKDcr1ZaScWM7VXZKeSl2u20yTgo0CkBtmuRcJSFdxqTW6r5Kr0hFY5imsxLl0`nType the 2nd string:`nConcatenated string is 'ZyBk27e3yEyMQ0WvL8t6KIJauYhrPmWuyRC2bnNKDcr1ZaScWM7VXZKeSl2u20yTgo0CkBtmuRcJSFdxqTW6r5Kr0hFY5imsxL0'`nType the 1st string (q - to quit):"
$ws.Range("F15").Value = "FAIL"
$ws.Range("G15").Value = "2nd string input is not being taken , validation should be added here for string 1. "

# Update the selected cell in the sheet view
$ws.Range("F4").Select()
